{"js": "// Applies the CV text-content edits described by the diff:\n//  1. \"employees outflow\"  -> \"employee's outflow\"\n//  2. \"Text-generating neural network (Torch-rnn, PyTorch, for social networks\n//     activities)\" -> \"Text-generating neural network (Tensorflow, Torch-rnn, PyTorch)\"\n//  3. \"...web-products (sites, Telegram-bots...\" -> \"...(web-sites, Telegram-bots...\"\n//  4. Technical skills: remove \"Google Tag Manager, \"\n//  5. Technical skills: append \", pytest, python-telegram-bot\" after \"bs4\"\n\nasync function findOne(context, searchText, options) {\n  const body = context.document.body;\n  const results = body.search(searchText, Object.assign({ matchCase: true }, options || {}));\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + searchText);\n  }\n  return results.items[0];\n}\n\nasync function replaceOnce(context, searchText, replacement, options) {\n  const range = await findOne(context, searchText, options);\n  range.insertText(replacement, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1. \"Predictive model of employees outflow\" -> \"...employee's outflow\"\nawait replaceOnce(context, \"employees\", \"employee's\");\n\n// 2a. Insert \"Tensorflow, \" right before the existing \"Torch-\" run.\nawait replaceOnce(context, \"Torch-\", \"Tensorflow, Torch-\");\n\n// 2b. Drop the trailing \", for social networks activities\" before the closing \")\".\nawait replaceOnce(context, \", for social networks activities)\", \")\");\n\n// 3. \"Creation and technical development of web-products (sites, Telegram-bots, ...)\"\n//    -> \"...(web-sites, Telegram-bots, ...)\"\nawait replaceOnce(context, \"(sites, Telegram-bots\", \"(web-sites, Telegram-bots\");\n\n// 4. Technical skills: drop \"Google Tag Manager, \"\nawait replaceOnce(context, \"Google Tag Manager, \", \"\");\n\n// 5. Technical skills: append \", pytest, python-telegram-bot\" right after \"bs4\"\n//    (insert after the existing range instead of replacing a span, so we don't\n//    disturb the surrounding proofErr markers).\nconst bs4Range = await findOne(context, \"bs4\");\nbs4Range.insertText(\", pytest, python-telegram-bot\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Applies the CV text-content edits described by the diff:\n#  1. \"employees outflow\"  -> \"employee's outflow\"\n#  2. \"Text-generating neural network (Torch-rnn, PyTorch, for social networks\n#     activities)\" -> \"Text-generating neural network (Tensorflow, Torch-rnn, PyTorch)\"\n#  3. \"...web-products (sites, Telegram-bots...\" -> \"...(web-sites, Telegram-bots...\"\n#  4. Technical skills: remove \"Google Tag Manager, \"\n#  5. Technical skills: append \", pytest, python-telegram-bot\" after \"bs4\"\n\n$d = $word.ActiveDocument\n\nfunction Find-And-SetText($searchText, $newText) {\n    $rng = $d.Content.Duplicate\n    $f = $rng.Find\n    $f.Text = $searchText\n    $f.MatchCase = $true\n    $f.MatchWholeWord = $false\n    $f.MatchWildcards = $false\n    $found = $f.Execute()\n    if (-not $found) {\n        throw \"Search text not found: $searchText\"\n    }\n    # Assign Range.Text directly (rather than Find.Replacement.Text) so no\n    # AutoCorrect/AutoFormat substitution (e.g. smart quotes) is applied.\n    $rng.Text = $newText\n    return $rng\n}\n\n# 1. \"Predictive model of employees outflow\" -> \"...employee's outflow\"\nFind-And-SetText \"employees\" \"employee's\" | Out-Null\n\n# 2a. Insert \"Tensorflow, \" right before the existing \"Torch-\" run.\nFind-And-SetText \"Torch-\" \"Tensorflow, Torch-\" | Out-Null\n\n# 2b. Drop the trailing \", for social networks activities\" before the closing \")\".\nFind-And-SetText \", for social networks activities)\" \")\" | Out-Null\n\n# 3. \"Creation and technical development of web-products (sites, Telegram-bots, ...)\"\n#    -> \"...(web-sites, Telegram-bots, ...)\"\nFind-And-SetText \"(sites, Telegram-bots\" \"(web-sites, Telegram-bots\" | Out-Null\n\n# 4. Technical skills: drop \"Google Tag Manager, \"\nFind-And-SetText \"Google Tag Manager, \" \"\" | Out-Null\n\n# 5. Technical skills: append \", pytest, python-telegram-bot\" right after \"bs4\"\n#    (insert after the found range instead of replacing it, so the surrounding\n#    proofErr markers around \"bs4\" are left untouched).\n$bs4Rng = $d.Content.Duplicate\n$f = $bs4Rng.Find\n$f.Text = \"bs4\"\n$f.MatchCase = $true\n$f.MatchWholeWord = $false\n$found = $f.Execute()\nif (-not $found) {\n    throw \"Search text not found: bs4\"\n}\n$bs4Rng.Collapse(0)  # wdCollapseEnd = 0\n$bs4Rng.InsertAfter(\", pytest, python-telegram-bot\")\n"}
